$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Numeric-looking text values (e.g. "1.0000") are forced to stay text
# by briefly applying a text NumberFormat, then restoring the default
# "Normal" cell style so the style index matches the original (unstyled) cell.

$ws.Cells.Item(2,4).Value = '25.970.76'
$ws.Cells.Item(2,5).Value = '  -0.52%  '
$ws.Cells.Item(3,4).Value = '1.745.80'
$ws.Cells.Item(3,5).Value = '  -0.10%  '
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '0.9999'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  +5.26%  '
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '1.0000'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  -0.03%  '
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.5073'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  -9.37%  '
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.2755'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = '  -2.70%  '
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06199'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  +0.32%  '
$ws.Cells.Item(10,4).Value = '1.746.65'
$ws.Cells.Item(10,5).Value = '  -0.46%  '
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07248'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  +0.73%  '
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.6542'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  +0.28%  '
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '15.19'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  -2.02%  '
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '4.673'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  +0.73%  '
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '77.79'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = '  -0.58%  '
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '1.000'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  +0.01%  '
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.9997'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(18,4).Value = '25.981.79'
$ws.Cells.Item(18,5).Value = '  -0.13%  '
$ws.Cells.Item(19,5).Value = '  +0.68%  '
$ws.Cells.Item(20,5).Value = '  +1.26%  '
$ws.Cells.Item(21,4).Value = '1.968.72'
$ws.Cells.Item(21,5).Value = '  -0.58%  '
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '4.431'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  +1.63%  '
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '8.749'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  -0.10%  '
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '5.393'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  +2.10%  '
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '136.71'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  -2.22%  '
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '1.518'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -0.65%  '
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '15.27'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  -0.55%  '
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.785'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = '  -1.36%  '
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '105.77'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  +0.39%  '
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '3.868'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  +1.83%  '
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.08202'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  -3.29%  '
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '3.655'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  +0.56%  '
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.04675'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  +0.46%  '
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '2.656'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  +0.28%  '
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '0.9986'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -0.89%  '
$ws.Cells.Item(36,5).Value = '  -1.99%  '
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '2.757'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  +1.88%  '
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.01617'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  +0.16%  '
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '1.933'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = '  -1.79%  '
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.9997'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  +0.02%  '
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '100.79'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  -0.28%  '
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.3927'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  -0.25%  '
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.7664'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = '  +2.57%  '
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '5.007'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = '  -1.21%  '
$ws.Cells.Item(45,5).Value = '  +0.27%  '
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '6.350'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  -0.06%  '
$ws.Cells.Item(47,2).Value = 'Aave'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '55.78'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  +1.80%  '
$ws.Cells.Item(48,2).Value = 'Cronos'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.05344'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = '  +0.13%  '
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '30.73'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  -0.24%  '
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.3452'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = '  -1.20%  '
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '7.596'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -0.30%  '
